# Weekly update: rotate the daily price-record data (columns D, J-Q) among the
# existing rows of the sheet. Row 3 keeps its original data; every other row
# (2, 4-15) receives the data that used to belong to a different row, per the
# mapping derived from the published diff.
#
# Columns that move together as one "record": D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# N (Unidad de comercializacion), O (Origen), P (Precio $/Kg), Q (Kg o Unidades)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# destinationRow -> sourceRow (source values, read from the ORIGINAL workbook
# state, are written into destination row)
$rowMap = @{
    2  = 13
    4  = 11
    5  = 9
    6  = 10
    7  = 6
    8  = 2
    9  = 7
    10 = 8
    11 = 12
    12 = 4
    13 = 14
    14 = 15
    15 = 5
}

# Snapshot the original values of every source row before writing anything,
# since several rows are both sources and destinations in this rotation.
$snapshot = @{}
foreach ($r in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($r)) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$r").Value2
        }
        $snapshot[$r] = $rowData
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowData[$col]
    }
}
